$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: zh-cn (E2) and de-de (F2) status cells
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-locale detail sheets: Status column (C2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the zh-cn / de-de status columns on the Overview sheet ---
$wsOverview.Columns("E").ColumnWidth = 12.5
$wsOverview.Columns("F").ColumnWidth = 12.5

# --- Narrow the Status column on each per-locale detail sheet ---
$wsZhCn.Columns("C").ColumnWidth = 12.5
$wsDeDe.Columns("C").ColumnWidth = 12.5
